# redo FR pop dens calcs using overlays
$wb = $excel.ActiveWorkbook

# area_mixre sheet: summary stats (std, 50%, 75%) recomputed
$wsMixre = $wb.Worksheets.Item("area_mixre")
$wsMixre.Range("B4").Value = 3.817123968925948
$wsMixre.Range("B7").Value = 4.373229608902978
$wsMixre.Range("B8").Value = 6.506528366714822

# area_pop_sum sheet: population + density recomputed
$wsPopSum = $wb.Worksheets.Item("area_pop_sum")
$wsPopSum.Range("B3").Value = 252120
$wsPopSum.Range("B4").Value = 1045.842797863891
